$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.286.48"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "3.523.16"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "584.21"
$ws.Range("E5").Value = "  +6.22%  "
$ws.Range("D6").Value = "179.27"
$ws.Range("E6").Value = "  -5.64%  "
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +3.96%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.640"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +7.01%  "
$ws.Range("D11").Value = "56.27"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Value = "'0.0000280"
$ws.Range("E12").Value = "  +4.27%  "
$ws.Range("D13").Value = "9.32"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "4.086.20"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "3.522.28"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "66.279.93"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "'12.10"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").Value = "416.29"
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("D22").Value = "4.31"
$ws.Range("E22").Value = "  +10.66%  "
$ws.Range("D23").Value = "4.44"
$ws.Range("E23").Value = "  +6.64%  "
$ws.Range("D24").Value = "85.69"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").Value = "13.36"
$ws.Range("E25").Value = "  +11.75%  "
$ws.Range("D26").Value = "11.15"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "6.05"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").Value = "9.13"
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("D30").Value = "30.48"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").Value = "'6.60"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").Value = "606.59"
$ws.Range("E32").Value = "  -6.48%  "
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "59.85"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("E36").Value = "  +9.57%  "
$ws.Range("D37").Value = "0.0₃0807"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "3.69"
$ws.Range("E38").Value = "  +11.37%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "37.34"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "3.270.11"
$ws.Range("E42").Value = "  +8.63%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "2.99"
$ws.Range("E44").Value = "  +4.22%  "
$ws.Range("D45").Value = "3.39"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("D47").Value = "0.0422"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "139.89"
$ws.Range("E51").Value = "  -0.01%  "
